$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 2.6
$ws.Range("Q2").Value = 1.95
$ws.Range("R2").Value = 1.9
$ws.Range("S2").Value = 2.6
$ws.Range("T2").Value = 1.48
$ws.Range("V2").Value = 1.23
$ws.Range("AF2").Value = 15
$ws.Range("AP2").Value = 17
$ws.Range("AR2").Value = 41
$ws.Range("G3").Value = 2.1
$ws.Range("H3").Value = 3.2
$ws.Range("I3").Value = 3.6
$ws.Range("J3").Value = 2.88
$ws.Range("K3").Value = 1.95
$ws.Range("L3").Value = 4.5
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("Q3").Value = 1.88
$ws.Range("R3").Value = 1.98
$ws.Range("S3").Value = 2.4
$ws.Range("T3").Value = 1.53
$ws.Range("U3").Value = 3.9
$ws.Range("AC3").Value = 6
$ws.Range("AD3").Value = 9
$ws.Range("AE3").Value = 9.5
$ws.Range("AF3").Value = 19
$ws.Range("AG3").Value = 21
$ws.Range("AJ3").Value = 6.5
$ws.Range("AN3").Value = 8
$ws.Range("AO3").Value = 17
$ws.Range("AP3").Value = 13
$ws.Range("AQ3").Value = 41
$ws.Range("AR3").Value = 34
$ws.Range("Q4").Value = 1.95
$ws.Range("R4").Value = 1.9
$ws.Range("S4").Value = 2.6
$ws.Range("T4").Value = 1.48
$ws.Range("G5").Value = 3.1
$ws.Range("H5").Value = 2.47
$ws.Range("I5").Value = 2.82
$ws.Range("J5").Value = 3.9
$ws.Range("L5").Value = 3.65
$ws.Range("M5").Value = 1.18
$ws.Range("N5").Value = 4.3
$ws.Range("O5").Value = 1.7
$ws.Range("P5").Value = 2.05
$ws.Range("S5").Value = 3.05
$ws.Range("T5").Value = 1.33
$ws.Range("W5").Value = 5.6
$ws.Range("X5").Value = 1.11
$ws.Range("Y5").Value = 1.72
$ws.Range("Z5").Value = 2.02
$ws.Range("AA5").Value = 2.25
$ws.Range("AB5").Value = 1.57
$ws.Range("AC5").Value = 6.2
$ws.Range("AD5").Value = 14
$ws.Range("AE5").Value = 12
$ws.Range("AF5").Value = 45
$ws.Range("AG5").Value = 40
$ws.Range("AH5").Value = 65
$ws.Range("AI5").Value = 4.3
$ws.Range("AJ5").Value = 5.2
$ws.Range("AK5").Value = 19
$ws.Range("AL5").Value = 150
$ws.Range("AN5").Value = 5.8
$ws.Range("AO5").Value = 12.5
$ws.Range("AP5").Value = 11.25
$ws.Range("AQ5").Value = 37
$ws.Range("AR5").Value = 35
$ws.Range("AS5").Value = 60
$ws.Range("G6").Value = 2.87
$ws.Range("H6").Value = 2.75
$ws.Range("I6").Value = 2.7
$ws.Range("J6").Value = 3.5
$ws.Range("K6").Value = 1.88
$ws.Range("L6").Value = 3.4
$ws.Range("M6").Value = 1.13
$ws.Range("N6").Value = 5.1
$ws.Range("O6").Value = 1.57
$ws.Range("P6").Value = 2.27
$ws.Range("S6").Value = 2.65
$ws.Range("T6").Value = 1.42
$ws.Range("W6").Value = 4.8
$ws.Range("X6").Value = 1.15
$ws.Range("Y6").Value = 1.55
$ws.Range("Z6").Value = 2.3
$ws.Range("AA6").Value = 2.1
$ws.Range("AB6").Value = 1.65
$ws.Range("AC6").Value = 6.6
$ws.Range("AD6").Value = 13
$ws.Range("AE6").Value = 11
$ws.Range("AF6").Value = 37
$ws.Range("AG6").Value = 32
$ws.Range("AH6").Value = 50
$ws.Range("AI6").Value = 5.1
$ws.Range("AJ6").Value = 5.5
$ws.Range("AK6").Value = 18
$ws.Range("AL6").Value = 120
$ws.Range("AM6").Value = 900
$ws.Range("AN6").Value = 6.2
$ws.Range("AO6").Value = 11.75
$ws.Range("AP6").Value = 10.75
$ws.Range("AQ6").Value = 32
$ws.Range("AR6").Value = 30
$ws.Range("AS6").Value = 50
$ws.Range("I7").Value = 3.3
$ws.Range("M7").Value = 1.13
$ws.Range("N7").Value = 6
$ws.Range("O7").Value = 1.62
$ws.Range("P7").Value = 2.3
$ws.Range("S7").Value = 3
$ws.Range("T7").Value = 1.4
$ws.Range("W7").Value = 6
$ws.Range("X7").Value = 1.13
$ws.Range("S8").Value = 2.1
$ws.Range("T8").Value = 1.7
$ws.Range("G9").Value = 2.6
$ws.Range("I9").Value = 2.55
$ws.Range("J9").Value = 3.4
$ws.Range("L9").Value = 3.4
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 10
$ws.Range("AC9").Value = 8
$ws.Range("AD9").Value = 12
$ws.Range("AE9").Value = 10
$ws.Range("AF9").Value = 26
$ws.Range("AG9").Value = 23
$ws.Range("AK9").Value = 15
$ws.Range("AM9").Value = 301
$ws.Range("AO9").Value = 12
$ws.Range("AP9").Value = 10
$ws.Range("AQ9").Value = 26
$ws.Range("AR9").Value = 21
$ws.Range("S10").Value = 2.15
$ws.Range("T10").Value = 1.67
$ws.Range("W10").Value = 4
$ws.Range("X10").Value = 1.22
$ws.Range("S11").Value = 2.2
$ws.Range("T11").Value = 1.65
$ws.Range("W11").Value = 4
$ws.Range("X11").Value = 1.22
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 3.3
$ws.Range("I14").Value = 3.5
$ws.Range("J14").Value = 2.62
$ws.Range("L14").Value = 4
$ws.Range("O14").Value = 1.3
$ws.Range("P14").Value = 3.4
$ws.Range("S14").Value = 2.03
$ws.Range("T14").Value = 1.83
$ws.Range("AB14").Value = 1.8
$ws.Range("AC14").Value = 7.5
$ws.Range("AD14").Value = 9.5
$ws.Range("AE14").Value = 9
$ws.Range("AF14").Value = 17
$ws.Range("AG14").Value = 17
$ws.Range("AI14").Value = 9.5
$ws.Range("AN14").Value = 11
$ws.Range("AO14").Value = 19
$ws.Range("AP14").Value = 13
$ws.Range("AQ14").Value = 41
$ws.Range("AR14").Value = 29
$ws.Range("AS14").Value = 41
